$d = $word.ActiveDocument

$pairs = @(
    @{old="224×5="; new="318×6="},
    @{old="823×9="; new="986×6="},
    @{old="389×5="; new="316×9="},
    @{old="548×6="; new="922×6="},
    @{old="512×4="; new="337×2="},
    @{old="892×4="; new="539×3="},
    @{old="605×9="; new="528×5="},
    @{old="251×2="; new="616×2="},
    @{old="523×6="; new="557×4="},
    @{old="930×5="; new="625×9="},
    @{old="591×4="; new="155×8="},
    @{old="662×9="; new="287×6="},
    @{old="536×6="; new="427×5="},
    @{old="156×7="; new="207×6="},
    @{old="123×4="; new="125×5="},
    @{old="303×6="; new="532×3="},
    @{old="396×6="; new="309×2="},
    @{old="927×9="; new="867×2="},
    @{old="221×8="; new="858×2="},
    @{old="145×5="; new="944×7="},
    @{old="928×7="; new="518×8="},
    @{old="549×7="; new="522×7="},
    @{old="707×6="; new="421×9="},
    @{old="142×8="; new="316×6="},
    @{old="708×4="; new="991×3="}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}
